$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.731.87'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '2.114.12'
$ws.Range("E3").Value = '  +10.19%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '255.44'
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.666'
$ws.Range("E6").Value = '  -4.14%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +6.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.94'
$ws.Range("E9").Value = '  +6.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.374'
$ws.Range("E10").Value = '  +2.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0742'
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.55'
$ws.Range("E13").Value = '  -0.20%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.415.46'
$ws.Range("E14").Value = '  +10.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.850'
$ws.Range("E15").Value = '  +6.63%  '
$ws.Range("D16").Value = '2.112.48'
$ws.Range("E16").Value = '  +10.28%  '
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '36.775.68'
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.35'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '0.0₃0847'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.49'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '241.69'
$ws.Range("E22").Value = '  -3.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.28'
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").Value = '  -8.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '172.68'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.35'
$ws.Range("E27").Value = '  +13.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.21'
$ws.Range("E28").Value = '  +4.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.03'
$ws.Range("E29").Value = '  -10.19%  '
$ws.Range("E30").Value = '  -4.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.38'
$ws.Range("E31").Value = '  +49.34%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0969'
$ws.Range("E33").Value = '  +13.84%  '
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").Value = '  +19.86%  '
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.19'
$ws.Range("E38").Value = '  -3.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.907'
$ws.Range("E39").Value = '  +4.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.36'
$ws.Range("E40").Value = '  -8.01%  '
$ws.Range("E41").Value = '  +6.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0223'
$ws.Range("E42").Value = '  -2.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.22'
$ws.Range("E43").Value = '  -4.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.85'
$ws.Range("E44").Value = '  +17.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.18'
$ws.Range("E45").Value = '  -5.54%  '
$ws.Range("D46").Value = '1.366.10'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0837'
$ws.Range("E47").Value = '  +3.48%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("E48").Value = '  -3.34%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.86'
$ws.Range("E49").Value = '  +6.78%  '
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.276.43'
$ws.Range("E51").Value = '  +9.23%  '
